# Generate Report for Archive
# - Update localization status text "Ready for handoff" -> "In Translation"
#   on all three sheets (Overview status columns E/F, and the Status column
#   on the zh-cn / de-de detail sheets).
# - Narrow the "Latest Handoff Datetime" / duplicate status-width columns
#   (Overview!E:F and Status col on zh-cn/de-de "C") to the new report width.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: columns E (zh-cn) and F (de-de), rows 2-3 ---
$wsOverview = $wb.Worksheets.Item("Overview")
$overviewRows = $wsOverview.UsedRange.Rows.Count
for ($r = 2; $r -le $overviewRows; $r++) {
    if ($wsOverview.Cells.Item($r, 5).Text -eq $oldStatus) {
        $wsOverview.Cells.Item($r, 5).Value = $newStatus
    }
    if ($wsOverview.Cells.Item($r, 6).Text -eq $oldStatus) {
        $wsOverview.Cells.Item($r, 6).Value = $newStatus
    }
}
# Resize the status columns to match the regenerated report layout.
$wsOverview.Range("E1").ColumnWidth = 12.5
$wsOverview.Range("F1").ColumnWidth = 12.5

# --- Per-locale detail sheets: "Status" column (C), rows 2-3 ---
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $ws.UsedRange.Rows.Count
    for ($r = 2; $r -le $rows; $r++) {
        if ($ws.Cells.Item($r, 3).Text -eq $oldStatus) {
            $ws.Cells.Item($r, 3).Value = $newStatus
        }
    }
    $ws.Range("C1").ColumnWidth = 12.5
}
